$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91-102 down to 92-103
$ws.Rows.Item(91).EntireRow.Insert()

# Populate the new row 91 with data (copy template columns from row 90, set specific values)
$ws.Cells.Item(91, 1).Value = 11
$ws.Cells.Item(91, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(91, 3).Value = "Bíobío"
$ws.Cells.Item(91, 4).Value = 44504
$ws.Cells.Item(91, 5).Value = 8
$ws.Cells.Item(91, 6).Value = 100112003
$ws.Cells.Item(91, 7).Value = "Ajo"
$ws.Cells.Item(91, 8).Value = "Chino"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 350
$ws.Cells.Item(91, 11).Value = 14000
$ws.Cells.Item(91, 12).Value = 15000
$ws.Cells.Item(91, 13).Value = 14571
$ws.Cells.Item(91, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(91, 15).Value = "China"
$ws.Cells.Item(91, 16).Value = 1457
$ws.Cells.Item(91, 17).Value = 10
$ws.Cells.Item(91, 18).Value = "Hortaliza"

# Ensure the date cell uses the same number format as the rest of column D
$ws.Cells.Item(91, 4).NumberFormat = $ws.Cells.Item(92, 4).NumberFormat
